{"js": "// Homework 1 final edit pass: extend the two \"Exercise 1.2\" answer\n// sentences with the reasoning clauses added in the commit\n// (\"o hw1 finish, new readings\").\n//\n//   1) \"... both a complexity of NP-hard. \" becomes\n//      \"... both a complexity of NP-hard, because they can both be\n//      reduced from the subset sum problem, which is NP-complete. \"\n//\n//   2) \"... time complexity O(n log n).\" becomes\n//      \"... time complexity O(n log n), because it iterates through\n//      all of the numbers in decreasing order.\"\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nparagraphs.items.forEach((p) => p.load(\"text\"));\nawait context.sync();\n\n// Locate the two paragraphs by a stable substring of their text rather\n// than a hard-coded index, so the script is resilient to unrelated\n// paragraph-count differences.\nlet npHardParagraph = null;\nlet bigOParagraph = null;\nfor (const p of paragraphs.items) {\n  if (p.text.indexOf(\"a complexity of NP-hard\") !== -1) {\n    npHardParagraph = p;\n  }\n  if (p.text.indexOf(\"time complexity O(n log n)\") !== -1) {\n    bigOParagraph = p;\n  }\n}\n\nif (!npHardParagraph) {\n  throw new Error(\"Could not find the 'complexity of NP-hard' paragraph.\");\n}\nif (!bigOParagraph) {\n  throw new Error(\"Could not find the 'time complexity O(n log n)' paragraph.\");\n}\n\n// --- Edit 1: \" NP-hard. \" -> \" NP-hard, because ... NP-complete. \"\nconst npHardHits = npHardParagraph.search(\"NP-hard. \", { matchCase: true });\nnpHardHits.load(\"items\");\nawait context.sync();\n\nif (npHardHits.items.length === 0) {\n  throw new Error(\"Could not find 'NP-hard. ' text to update.\");\n}\nnpHardHits.items[0].insertText(\n  \"NP-hard, because they can both be reduced from the subset sum problem, which is NP-complete. \",\n  Word.InsertLocation.replace\n);\nawait context.sync();\n\n// --- Edit 2: \"O(n log n).\" -> \"O(n log n), because ... decreasing order.\"\nconst bigOHits = bigOParagraph.search(\"O(n log n).\", { matchCase: true });\nbigOHits.load(\"items\");\nawait context.sync();\n\nif (bigOHits.items.length === 0) {\n  throw new Error(\"Could not find 'O(n log n).' text to update.\");\n}\nbigOHits.items[0].insertText(\n  \"O(n log n), because it iterates through all of the numbers in decreasing order.\",\n  Word.InsertLocation.replace\n);\nawait context.sync();\n", "ps1": "# Homework 1 final edit pass: extend the two \"Exercise 1.2\" answer\n# sentences with the reasoning clauses added in the commit\n# (\"o hw1 finish, new readings\").\n#\n#   1) \"... both a complexity of NP-hard. \" becomes\n#      \"... both a complexity of NP-hard, because they can both be\n#      reduced from the subset sum problem, which is NP-complete. \"\n#\n#   2) \"... time complexity O(n log n).\" becomes\n#      \"... time complexity O(n log n), because it iterates through\n#      all of the numbers in decreasing order.\"\n\n$d = $word.ActiveDocument\n\n# Locate the two paragraphs by a stable substring of their text rather\n# than a hard-coded index, so the script is resilient to unrelated\n# paragraph-count differences.\n$npHardParagraph = $null\n$bigOParagraph = $null\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $t = $p.Range.Text\n    if ($t -like \"*a complexity of NP-hard*\") {\n        $npHardParagraph = $p\n    }\n    if ($t -like \"*time complexity O(n log n)*\") {\n        $bigOParagraph = $p\n    }\n}\n\nif ($null -eq $npHardParagraph) {\n    throw \"Could not find the 'complexity of NP-hard' paragraph.\"\n}\nif ($null -eq $bigOParagraph) {\n    throw \"Could not find the 'time complexity O(n log n)' paragraph.\"\n}\n\n# --- Edit 1: \" NP-hard. \" -> \" NP-hard, because ... NP-complete. \"\n$r1 = $npHardParagraph.Range\n$r1.Find.ClearFormatting()\n$r1.Find.Text = \"NP-hard. \"\n$r1.Find.MatchCase = $true\n$r1.Find.MatchWholeWord = $false\n$r1.Find.Forward = $true\n$r1.Find.Wrap = 0\n$found1 = $r1.Find.Execute()\nif (-not $found1) {\n    throw \"Could not find 'NP-hard. ' text to update.\"\n}\n$r1.Text = \"NP-hard, because they can both be reduced from the subset sum problem, which is NP-complete. \"\n\n# --- Edit 2: \"O(n log n).\" -> \"O(n log n), because ... decreasing order.\"\n$r2 = $bigOParagraph.Range\n$r2.Find.ClearFormatting()\n$r2.Find.Text = \"O(n log n).\"\n$r2.Find.MatchCase = $true\n$r2.Find.MatchWholeWord = $false\n$r2.Find.Forward = $true\n$r2.Find.Wrap = 0\n$found2 = $r2.Find.Execute()\nif (-not $found2) {\n    throw \"Could not find 'O(n log n).' text to update.\"\n}\n$r2.Text = \"O(n log n), because it iterates through all of the numbers in decreasing order.\"\n"}
